$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览 = Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 63
$ws1.Range("F15").Value = 54
$ws1.Range("F18").Value = 593
$ws1.Range("F20").Value = 5784
$ws1.Range("F26").Value = 5443
$ws1.Range("F27").Value = 5443
$ws1.Range("F30").Value = 1561
$ws1.Range("F31").Value = 398
$ws1.Range("F34").Value = 1092
$ws1.Range("F35").Value = 673
$ws1.Range("F36").Value = 137
$ws1.Range("F38").Value = 70

# Sheet "演出" (演出 = Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 70

# Sheet "本地生活" (本地生活 = Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 285

# Sheet "全部类型" (全部类型 = All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 63
$ws4.Range("F15").Value = 54
$ws4.Range("F19").Value = 593
$ws4.Range("F21").Value = 5784
$ws4.Range("F29").Value = 5443
$ws4.Range("F30").Value = 5443
$ws4.Range("F33").Value = 1561
$ws4.Range("F34").Value = 399
$ws4.Range("F36").Value = 1092
$ws4.Range("F37").Value = 673
$ws4.Range("F38").Value = 137
$ws4.Range("F44").Value = 70

$wb.Save()
